# HIVE TEAMS.docx -- English to Russian translation pass
#
# Applies the text substitutions described by the commit diff. Uses
# Find.Execute purely to *locate* each target range (no in-place find
# replacement text, since Find's own replace path silently "smart
# quotes"/autocorrects the replacement string); the actual text swap is
# then done with a plain Range.Text assignment, which leaves the
# characters exactly as supplied. Longer / more specific phrases are
# replaced before shorter substrings they contain, so ordering matters
# (e.g. "Front-End Developer" before bare "Developer").

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# Replace the single (first) occurrence of $find with $repl.
function Replace-One($find, $repl) {
    $rng = $d.Content
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if ($found) {
        $rng.Text = $repl
    }
    return $found
}

# Replace every occurrence of $find with $repl.
function Replace-All($find, $repl) {
    $count = 0
    while ($true) {
        $rng = $d.Content
        $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
        if (-not $found) { break }
        $rng.Text = $repl
        $count = $count + 1
        if ($count -gt 50) { break }
    }
    return $count
}

# --- HIVE TEAM: WEB -----------------------------------------------------
Replace-One ($nbsp + "HIVE TEAM: WEB") "Команда Hive: Web"
Replace-One "Responsible for building web applications." "Создание веб-приложений и различных инструментов."
Replace-All "Hive Coordinator" "Координатор Hive"
Replace-One "From the Guetos of Brazil to the cryptoverse. Father, husband and technolover, born in Brazil and living in Australia. Mobile and Web Development." "От гетто в Бразилии до мира криптовалют. Отец, муж и настоящий энтузиаст, увлечённый технологиями. Родился в Бразилии, но сейчас живёт в Австралии. Занимается мобильной и веб-разработкой."
Replace-One "Front-End Developer" "Front-End Разработчик"
Replace-All "Developer" "Разработчик"

# --- HIVE TEAM: QUALITY ASSURANCE (two runs merge into one) -------------
$rng = $d.Content
$rng.Find.Execute($nbsp + "HIVE TEAM: QUALITY ASSURANCE", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rng.Text = "Команда Hive: Контроль Качества"
$rng.Font.Name = "Helvetica Neue"
$rng.Font.Bold = $false
$rng.Font.BoldBi = $false

Replace-One "Ensuring all development tasks meet quality criteria." "Тестирование и контроль качества на всех этапах разработки."
Replace-One "Release Coordinator" "Релиз-координатор"
Replace-One "auditor" "Аудитор"
Replace-One "Senior QA Tester" "Специалист по обеспечению качества"
Replace-One "security" "Безопасность"
Replace-One "QA Tester" "Тестировщик"

# --- HIVE TEAM: DEVELOPMENT ----------------------------------------------
Replace-One ($nbsp + "HIVE TEAM: DEVELOPMENT") "Команда Hive: Разработка"
Replace-One "Responsible for building SmartCash and supporting applications." "Разработка экосистемы SmartCash, сервисы и приложения."
Replace-One "Creator of the Dash N Drink Soda Machine & SmartCash POS." "Создатель Dash N Drink Soda Machine и SmartCash POS."
Replace-All "Developer" "Разработчик"

# --- HIVE TEAM: OUTREACH 2 ------------------------------------------------
Replace-One ($nbsp + "HIVE TEAM: OUTREACH 2") "Команда Hive: Продвижение (II)"
Replace-One "This team focuses on community building, growth, general user acquisition in South America" "Работа над расширением сообщества в Латинской Америке."
Replace-All "Outreach Support" "Продвижение"

# --- HIVE TEAM: SUPPORT & WEB ---------------------------------------------
Replace-One ($nbsp + "HIVE TEAM: SUPPORT" + $nbsp) "Команда Hive: Поддержка "
Replace-One ($nbsp + "WEB") "Web"
Replace-One "This Hive is responsible for on-boarding & generalized SmartCash support." "Интеграции и поддержка пользователей."
Replace-One "Alex is a jack of all trades who loves Technology, Graphics, Web Design & Infrastructure." 'Alex – настоящий "человек всех профессий". Его конёк - технологии, графика, работа с инфраструктурой сайтов и веб-дизайн.'
Replace-One "Fiscal Officer" "Финансовый управляющий"
Replace-All "Support" "Поддержка"
Replace-One "Vice Coordinator" "Вице-координатор"
Replace-One "Video Guru" "Видео-гуру"
Replace-One "Legal affairs" "Юридические вопросы"

# --- Footer / closing section ---------------------------------------------
Replace-One "WANNA GET INVOLVED?" "Хотите присоединиться?"
Replace-One "The SmartHive has a place for folks of all backgrounds. Come hungry!" "SmartHive – это место, где ваши таланты и способности будут востребованы. Присоединяйтесь!  "
Replace-One "We believe ‘Core’ teams are a bad idea and something that ultimately leads to inefficiency and corruption. We want to move past it and create a decentralized organizational model inspired by ant and bee colonies." "Мы считаем, что постоянные команды — это то, что приводит к коррупции и неэффективности. Мы хотим идти по собственному пути и поэтому создали децентрализованную организационную модель, основанную на принципах жизни муравьёв и пчелиных колоний."
Replace-One "In order to create and maintain a decentralized governance structure, we are introducing two concepts SmartHive and Hive Structuring Teams (HST). SmartHive enables anyone that holds coins the opportunity to vote on proposals submitted by the community. SmartHive will be the lifeblood of the project, which will allow anyone to get involved and submit proposals, helping to generate organic growth at a grassroots level, creating a bottom-up management structure." "Создание и поддержание такой структуры управления требует особого подхода, поэтому мы разработали две концепции – SmartHive и Hive Structuring Teams (HST).  SmartHive дает возможность любому держателю монет голосовать за проекты и идеи, представленные сообществом.  Благодаря SmartHive каждый участник способен проявлять себя – выдвигать свои идеи на голосование, участвовать в обсуждениях, способствовать росту сообщества, а также голосовать за другие предложения."
